# Auto-generated Excel COM-interop script that applies the numeric
# corrections described by the commit diff to the "Cactuar_Profits"
# workbook. The workbook contains 8 sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR); each sheet has the same "Leve" profit-tracking layout
# in columns A:N. Only specific H..N cells on specific rows changed;
# everything else in the workbook is left untouched.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 1000
$ws.Range("N19").Value = -1350
# Row 40
$ws.Range("H40").Value = 17879414
$ws.Range("I40").Value = 9836.4
$ws.Range("J40").Value = 38498156
$ws.Range("K40").Value = 9836.4
$ws.Range("L40").Value = 38498156
$ws.Range("M40").Value = -9661.4
$ws.Range("N40").Value = -38498506
# Row 58
$ws.Range("H58").Value = 3887
$ws.Range("I58").Value = 274.5
$ws.Range("K58").Value = 823.5
$ws.Range("M58").Value = -673.5
# Row 62
$ws.Range("H62").Value = 3573
$ws.Range("I62").Value = 3201
$ws.Range("J62").Value = 3883
$ws.Range("K62").Value = 3201
$ws.Range("L62").Value = 3883
$ws.Range("M62").Value = -2577
$ws.Range("N62").Value = -5131
# Row 65
$ws.Range("H65").Value = 3573
$ws.Range("I65").Value = 3201
$ws.Range("J65").Value = 3883
$ws.Range("K65").Value = 16005
$ws.Range("L65").Value = 19415
$ws.Range("M65").Value = -12885
$ws.Range("N65").Value = -25655
# Row 86
$ws.Range("H86").Value = 1248453
$ws.Range("I86").Value = 2160193.8
$ws.Range("K86").Value = 2160193.8
$ws.Range("M86").Value = -2159070.8
# Row 89
$ws.Range("H89").Value = 1248453
$ws.Range("I89").Value = 2160193.8
$ws.Range("K89").Value = 10800969
$ws.Range("M89").Value = -10795353
# Row 98
$ws.Range("H98").Value = 2898.4443
$ws.Range("I98").Value = 1849.3334
$ws.Range("K98").Value = 1849.3334
$ws.Range("M98").Value = -351.3334
# Row 122
$ws.Range("H122").Value = 2898.4443
$ws.Range("I122").Value = 1849.3334
$ws.Range("K122").Value = 5548.0002
$ws.Range("M122").Value = -3098.0002
# Row 125
$ws.Range("H125").Value = 2450.2
$ws.Range("I125").Value = 2031
$ws.Range("J125").Value = 2514.6924
$ws.Range("K125").Value = 18279
$ws.Range("L125").Value = 22632.2316
$ws.Range("M125").Value = -15819
$ws.Range("N125").Value = -27552.2316
# Row 138
$ws.Range("H138").Value = 4676.78
$ws.Range("I138").Value = 1790.2858
$ws.Range("K138").Value = 5370.857400000001
$ws.Range("M138").Value = -230.8574000000008
# Row 139
$ws.Range("H139").Value = 125000
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280
# Row 140
$ws.Range("H140").Value = 68673.11
$ws.Range("J140").Value = 68418.625
$ws.Range("L140").Value = 68418.625
$ws.Range("N140").Value = -78778.625
# Remove cell M19 entirely (no longer present after edit)
$ws.Range("M19").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 35062.2
$ws.Range("I61").Value = 53103.668
$ws.Range("K61").Value = 53103.668
$ws.Range("M61").Value = -52891.668
# Row 110
$ws.Range("H110").Value = 3408731
$ws.Range("I110").Value = 5103596.5
$ws.Range("K110").Value = 5103596.5
$ws.Range("M110").Value = -5101551.5
# Row 128
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
# Row 136
$ws.Range("H136").Value = 35062.2
$ws.Range("I136").Value = 53103.668
$ws.Range("K136").Value = 159311.004
$ws.Range("M136").Value = -156761.004

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2583.0908
$ws.Range("I20").Value = 3046.1538
$ws.Range("J20").Value = 1914.2222
$ws.Range("K20").Value = 3046.1538
$ws.Range("L20").Value = 1914.2222
$ws.Range("M20").Value = -2799.1538
$ws.Range("N20").Value = -2408.2222
# Row 134
$ws.Range("H134").Value = 1593.1111
$ws.Range("I134").Value = 1593.1111
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4779.3333
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2244.3333
# Remove cell N134 entirely (no longer present after edit)
$ws.Range("N134").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2033.8667
$ws.Range("I16").Value = 2222.2
$ws.Range("J16").Value = 1939.7
$ws.Range("K16").Value = 2222.2
$ws.Range("L16").Value = 1939.7
$ws.Range("M16").Value = -1935.2
$ws.Range("N16").Value = -2513.7
# Row 58
$ws.Range("H58").Value = 457480.97
$ws.Range("I58").Value = 1001980.7
$ws.Range("J58").Value = 3731.1667
$ws.Range("K58").Value = 1001980.7
$ws.Range("L58").Value = 3731.1667
$ws.Range("M58").Value = -1001777.7
$ws.Range("N58").Value = -4137.1667
# Row 62
$ws.Range("H62").Value = 46953
$ws.Range("I62").Value = 4002.5
$ws.Range("K62").Value = 4002.5
$ws.Range("M62").Value = -3378.5
# Row 65
$ws.Range("H65").Value = 46953
$ws.Range("I65").Value = 4002.5
$ws.Range("K65").Value = 20012.5
$ws.Range("M65").Value = -16892.5
# Row 99
$ws.Range("H99").Value = 13374.167
$ws.Range("I99").Value = 30485
$ws.Range("J99").Value = 8485.357
$ws.Range("K99").Value = 30485
$ws.Range("L99").Value = 8485.357
$ws.Range("M99").Value = -28987
$ws.Range("N99").Value = -11481.357
# Row 105
$ws.Range("H105").Value = 1337806.1
$ws.Range("I105").Value = 1748923.4
$ws.Range("K105").Value = 1748923.4
$ws.Range("M105").Value = -1747176.4
# Row 113
$ws.Range("H113").Value = 2033.8667
$ws.Range("I113").Value = 2222.2
$ws.Range("J113").Value = 1939.7
$ws.Range("K113").Value = 2222.2
$ws.Range("L113").Value = 1939.7
$ws.Range("M113").Value = -52.19999999999982
$ws.Range("N113").Value = -6279.7
# Row 126
$ws.Range("H126").Value = 13374.167
$ws.Range("I126").Value = 30485
$ws.Range("J126").Value = 8485.357
$ws.Range("K126").Value = 91455
$ws.Range("L126").Value = 25456.071
$ws.Range("M126").Value = -88985
$ws.Range("N126").Value = -30396.071
# Row 132
$ws.Range("H132").Value = 16679096
$ws.Range("I132").Value = 22238034
$ws.Range("J132").Value = 2280
$ws.Range("K132").Value = 66714102
$ws.Range("L132").Value = 6840
$ws.Range("M132").Value = -66711572
$ws.Range("N132").Value = -11900
# Row 134
$ws.Range("H134").Value = 4263.857
$ws.Range("I134").Value = 4169.4
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 12508.2
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -9973.199999999999
$ws.Range("N134").Value = -18570
# Row 136
$ws.Range("H136").Value = 457480.97
$ws.Range("I136").Value = 1001980.7
$ws.Range("J136").Value = 3731.1667
$ws.Range("K136").Value = 3005942.1
$ws.Range("L136").Value = 11193.5001
$ws.Range("M136").Value = -3003392.1
$ws.Range("N136").Value = -16293.5001
# Row 138
$ws.Range("H138").Value = 85610.336
$ws.Range("J138").Value = 81499.14
$ws.Range("L138").Value = 81499.14
$ws.Range("N138").Value = -91779.14
# Row 139
$ws.Range("H139").Value = 100499.25
$ws.Range("J139").Value = 100499.25
$ws.Range("L139").Value = 100499.25
$ws.Range("N139").Value = -110779.25
# Row 141
$ws.Range("H141").Value = 79208.336
$ws.Range("J141").Value = 80948.8
$ws.Range("L141").Value = 80948.8
$ws.Range("N141").Value = -91308.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 134
$ws.Range("H134").Value = 8705.105
$ws.Range("I134").Value = 2363.4666
$ws.Range("K134").Value = 7090.399800000001
$ws.Range("M134").Value = -2020.399800000001

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 2272506.5
$ws.Range("I70").Value = 3406360.8
$ws.Range("K70").Value = 3406360.8
$ws.Range("M70").Value = -3406090.8
# Row 73
$ws.Range("H73").Value = 2272506.5
$ws.Range("I73").Value = 3406360.8
$ws.Range("K73").Value = 3406360.8
$ws.Range("M73").Value = -3405424.8
# Row 94
$ws.Range("H94").Value = 45149.332
$ws.Range("J94").Value = 35224
$ws.Range("L94").Value = 35224
$ws.Range("N94").Value = -36576

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 46155776
$ws.Range("I16").Value = 63159276
$ws.Range("J16").Value = 3414.4285
$ws.Range("K16").Value = 63159276
$ws.Range("L16").Value = 3414.4285
$ws.Range("M16").Value = -63159106
$ws.Range("N16").Value = -3754.4285
# Row 46
$ws.Range("H46").Value = 5864.816
$ws.Range("J46").Value = 5919
$ws.Range("L46").Value = 5919
$ws.Range("N46").Value = -6295
# Row 68
$ws.Range("H68").Value = 2526992
$ws.Range("I68").Value = 2842449.5
$ws.Range("J68").Value = 3332
$ws.Range("K68").Value = 2842449.5
$ws.Range("L68").Value = 3332
$ws.Range("M68").Value = -2841700.5
$ws.Range("N68").Value = -4830
# Row 71
$ws.Range("H71").Value = 2526992
$ws.Range("I71").Value = 2842449.5
$ws.Range("J71").Value = 3332
$ws.Range("K71").Value = 14212247.5
$ws.Range("L71").Value = 16660
$ws.Range("M71").Value = -14208503.5
$ws.Range("N71").Value = -24148
# Row 82
$ws.Range("H82").Value = 5209033
$ws.Range("I82").Value = 7812687
$ws.Range("K82").Value = 7812687
$ws.Range("M82").Value = -7812326
# Row 85
$ws.Range("H85").Value = 5209033
$ws.Range("I85").Value = 7812687
$ws.Range("K85").Value = 7812687
$ws.Range("M85").Value = -7811439

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 37045068
$ws.Range("I132").Value = 6945631
$ws.Range("K132").Value = 20836893
$ws.Range("M132").Value = -20834363
# Row 133
$ws.Range("H133").Value = 67994.5
$ws.Range("J133").Value = 67994.5
$ws.Range("L133").Value = 67994.5
$ws.Range("N133").Value = -78114.5

